# Rename book to story (#87)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C5").Value = "Story.csv"
$ws.Range("C6").Value = "StoryChapter.csv"

$ws.Range("C23").Select()
